$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.4523597029431983
$ws.Range("C2").Value = 0.1287589183487867
$ws.Range("E2").Value = 0.1171229856424318
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002481159070970639
$ws.Range("I2").Value = 0.9845259143874507
$ws.Range("K2").Value = 0.2757525494700133
$ws.Range("L2").Value = 0.2112457154521081
$ws.Range("M2").Value = 0.1439763166999981
$ws.Range("O2").Value = 3.726952385020709
# Row 3
$ws.Range("B3").Value = 0.4162846600611942
$ws.Range("C3").Value = 0.1276140261129797
$ws.Range("E3").Value = 0.1176478481342205
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002483328078575043
$ws.Range("I3").Value = 0.9956685842478841
$ws.Range("K3").Value = 0.2429546270339245
$ws.Range("L3").Value = 0.208750502066934
$ws.Range("M3").Value = 0.1370848447005635
$ws.Range("O3").Value = 3.766878640232349
# Row 4
$ws.Range("B4").Value = 0.3941993028792865
$ws.Range("C4").Value = 0.1269118651315111
$ws.Range("E4").Value = 0.1180178252732205
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.002484730859651526
$ws.Range("I4").Value = 1.002976018568582
$ws.Range("K4").Value = 0.2227811129830997
$ws.Range("L4").Value = 0.2073195973426252
$ws.Range("M4").Value = 0.1329005115409032
$ws.Range("O4").Value = 3.793306284234376
# Row 5
$ws.Range("B5").Value = 0.3852163118349381
$ws.Range("C5").Value = 0.1266259592590373
$ws.Range("E5").Value = 0.1181806122763227
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002485320410211591
$ws.Range("I5").Value = 1.006071009493059
$ws.Range("K5").Value = 0.2145518526503736
$ws.Range("L5").Value = 0.2067620058614992
$ws.Range("M5").Value = 0.1312073283952166
$ws.Range("O5").Value = 3.804556776888262
# Row 6
$ws.Range("B6").Value = 0.3837257375166985
$ws.Range("C6").Value = 0.1265784995356185
$ws.Range("E6").Value = 0.1182083694391665
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002485419387507209
$ws.Range("I6").Value = 1.00659200682796
$ws.Range("K6").Value = 0.2131848990620284
$ws.Range("L6").Value = 0.2066709613903939
$ws.Range("M6").Value = 0.1309269033057561
$ws.Range("O6").Value = 3.806453965636848
# Row 7
$ws.Range("B7").Value = 0.3940780855270702
$ws.Range("C7").Value = 0.12690800833537
$ws.Range("E7").Value = 0.1180199719854613
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002484738737989726
$ws.Range("I7").Value = 1.003017284320432
$ws.Range("K7").Value = 0.2226701635703563
$ws.Range("L7").Value = 0.207311974065
$ws.Range("M7").Value = 0.1328776280514425
$ws.Range("O7").Value = 3.793456064665079
# Row 8
$ws.Range("B8").Value = 0.4399079496160709
$ws.Range("C8").Value = 0.1283640079237145
$ws.Range("E8").Value = 0.1172940688697981
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002481892238955574
$ws.Range("I8").Value = 0.9882712483477931
$ws.Range("K8").Value = 0.2644514902856656
$ws.Range("L8").Value = 0.2103644123352879
$ws.Range("M8").Value = 0.1415904477198353
$ws.Range("O8").Value = 3.740321986186132
# Row 9
$ws.Range("B9").Value = 0.5302704066480715
$ws.Range("C9").Value = 0.1312245696748917
$ws.Range("E9").Value = 0.1162482338601052
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002476871280257576
$ws.Range("I9").Value = 0.9630487987687779
$ws.Range("K9").Value = 0.3460840320424268
$ws.Range("L9").Value = 0.2171501966113283
$ws.Range("M9").Value = 0.159044739797217
$ws.Range("O9").Value = 3.651301402494781
# Row 10
$ws.Range("B10").Value = 0.5969324507784393
$ws.Range("C10").Value = 0.1333282834536931
$ws.Range("E10").Value = 0.1157089544541385
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002473521088188213
$ws.Range("I10").Value = 0.9467679581041963
$ws.Range("K10").Value = 0.4058563027639366
$ws.Range("L10").Value = 0.2226204855192577
$ws.Range("M10").Value = 0.1720878438163709
$ws.Range("O10").Value = 3.595146510079303
# Row 11
$ws.Range("B11").Value = 0.6273128576712566
$ws.Range("C11").Value = 0.1342855251342669
$ws.Range("E11").Value = 0.1155131324669121
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002472069842852375
$ws.Range("I11").Value = 0.9398495092083685
$ws.Range("K11").Value = 0.4330002813868248
$ws.Range("L11").Value = 0.2252137574048874
$ws.Range("M11").Value = 0.1780681226675398
$ws.Range("O11").Value = 3.571607914037386
# Row 12
$ws.Range("B12").Value = 0.6388245321844295
$ws.Range("C12").Value = 0.1346480194514683
$ws.Range("E12").Value = 0.1154460772779391
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.00247153070588587
$ws.Range("I12").Value = 0.9372997948468402
$ws.Range("K12").Value = 0.4432718265495055
$ws.Range("L12").Value = 0.2262107690738304
$ws.Range("M12").Value = 0.1803393275260774
$ws.Range("O12").Value = 3.562983000017724
# Row 13
$ws.Range("B13").Value = 0.6363449748627374
$ws.Range("C13").Value = 0.134569949926636
$ws.Range("E13").Value = 0.1154602034229715
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002471646356010548
$ws.Range("I13").Value = 0.9378458012787902
$ws.Range("K13").Value = 0.4410599988362662
$ws.Range("L13").Value = 0.2259953790321276
$ws.Range("M13").Value = 0.1798498915076436
$ws.Range("O13").Value = 3.564827690534131
# Row 14
$ws.Range("B14").Value = 0.6282597883351286
$ws.Range("C14").Value = 0.1343153477688261
$ws.Range("E14").Value = 0.1155074736466677
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002472025279103355
$ws.Range("I14").Value = 0.9396383366694785
$ws.Range("K14").Value = 0.4338454769655868
$ws.Range("L14").Value = 0.2252954820490345
$ws.Range("M14").Value = 0.1782548445353314
$ws.Range("O14").Value = 3.570892551778769
# Row 15
$ws.Range("B15").Value = 0.6233083044029399
$ws.Range("C15").Value = 0.134159396752743
$ws.Range("E15").Value = 0.1155373518580713
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002472258735851478
$ws.Range("I15").Value = 0.940745452493907
$ws.Range("K15").Value = 0.4294254076539517
$ws.Range("L15").Value = 0.2248687255964654
$ws.Range("M15").Value = 0.1772786879125121
$ws.Range("O15").Value = 3.574645047089334
# Row 16
$ws.Range("B16").Value = 0.5949480761209713
$ws.Range("C16").Value = 0.1332657280795644
$ws.Range("E16").Value = 0.1157227462564254
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002473617390879345
$ws.Range("I16").Value = 0.9472299112902505
$ws.Range("K16").Value = 0.4040813907911911
$ws.Range("L16").Value = 0.2224531116835777
$ws.Range("M16").Value = 0.1716979502516907
$ws.Range("O16").Value = 3.596725198294678
# Row 17
$ws.Range("B17").Value = 0.5775637071417918
$ws.Range("C17").Value = 0.1327175354032732
$ws.Range("E17").Value = 0.1158491436757085
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.00247446948758627
$ws.Range("I17").Value = 0.9513328569864008
$ws.Range("K17").Value = 0.3885212824527855
$ws.Range("L17").Value = 0.2209980013556532
$ws.Range("M17").Value = 0.1682862610313123
$ws.Range("O17").Value = 3.610784653478419
# Row 18
$ws.Range("B18").Value = 0.5675699448531759
$ws.Range("C18").Value = 0.1324022548610913
$ws.Range("E18").Value = 0.1159265043694653
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002474966443615853
$ws.Range("I18").Value = 0.9537386751254147
$ws.Range("K18").Value = 0.3795671555319586
$ws.Range("L18").Value = 0.2201709301299672
$ws.Range("M18").Value = 0.1663283710560677
$ws.Range("O18").Value = 3.619060139729555
# Row 19
$ws.Range("B19").Value = 0.5641871548785389
$ws.Range("C19").Value = 0.1322955115186275
$ws.Range("E19").Value = 0.1159534983671886
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002475135882768567
$ws.Range("I19").Value = 0.9545611290240466
$ws.Range("K19").Value = 0.3765347098944289
$ws.Range("L19").Value = 0.2198925954658222
$ws.Range("M19").Value = 0.1656662272067067
$ws.Range("O19").Value = 3.621894510860614
# Row 20
$ws.Range("B20").Value = 0.5794137629237639
$ws.Range("C20").Value = 0.1327758890284585
$ws.Range("E20").Value = 0.1158352062788097
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002474378071531749
$ws.Range("I20").Value = 0.9508913395985807
$ws.Range("K20").Value = 0.3901781383653997
$ws.Range("L20").Value = 0.2211518794003524
$ws.Range("M20").Value = 0.168648984393414
$ws.Range("O20").Value = 3.609268453479217
# Row 21
$ws.Range("B21").Value = 0.6306344118330287
$ws.Range("C21").Value = 0.1343901306124025
$ws.Range("E21").Value = 0.1154933967446894
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002471913697922558
$ws.Range("I21").Value = 0.9391099215943193
$ws.Range("K21").Value = 0.4359647593248042
$ws.Range("L21").Value = 0.2255006523784147
$ws.Range("M21").Value = 0.1787231704206462
$ws.Range("O21").Value = 3.569103320553765
# Row 22
$ws.Range("B22").Value = 0.6641521366140637
$ws.Range("C22").Value = 0.1354451624029167
$ws.Range("E22").Value = 0.1153113699881985
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002470363791065726
$ws.Range("I22").Value = 0.9318189881689989
$ws.Range("K22").Value = 0.4658462049028174
$ws.Range("L22").Value = 0.2284301992072102
$ws.Range("M22").Value = 0.1853456552612371
$ws.Range("O22").Value = 3.544535555111921
# Row 23
$ws.Range("B23").Value = 0.6462594802314072
$ws.Range("C23").Value = 0.134882078562903
$ws.Range("E23").Value = 0.1154047426005569
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002471185466395473
$ws.Range("I23").Value = 0.9356728787972273
$ws.Range("K23").Value = 0.4499020156420954
$ws.Range("L23").Value = 0.2268586753963291
$ws.Range("M23").Value = 0.1818076416906749
$ws.Range("O23").Value = 3.55749385383821
# Row 24
$ws.Range("B24").Value = 0.578577350337099
$ws.Range("C24").Value = 0.1327495077170582
$ws.Range("E24").Value = 0.1158414927560916
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002474419378642296
$ws.Range("I24").Value = 0.9510908032707412
$ws.Range("K24").Value = 0.3894290999713803
$ws.Range("L24").Value = 0.2210822815673481
$ws.Range("M24").Value = 0.1684849861379547
$ws.Range("O24").Value = 3.609953327824797
# Row 25
$ws.Range("B25").Value = 0.5057751927134291
$ws.Range("C25").Value = 0.1304502642380356
$ws.Range("E25").Value = 0.1164908463546901
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002478169865152813
$ws.Range("I25").Value = 0.9694768446622497
$ws.Range("K25").Value = 0.3240345090830772
$ws.Range("L25").Value = 0.2152291207016077
$ws.Range("M25").Value = 0.1542839815318438
$ws.Range("O25").Value = 3.673759499059074
